# TC_164_Verify_DC_WorstCase_Units_Indicator.xlsx
# Updated test data for DC, TripCurrent, Voltdrop, BatteryStandby
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add new "Loop" / "Column" section in columns H:I -----------------
# Reuse the existing header accent style (currently on F1) for the two
# new header cells, then overwrite their text.
$ws.Range("F1").Copy($ws.Range("H1"))
$ws.Range("F1").Copy($ws.Range("I1"))
$ws.Range("H1").Value = "Loop"
$ws.Range("I1").Value = "Column"

# Reuse the plain bordered style (currently on B4) for the new data
# cells, then overwrite their text/values.
$ws.Range("B4").Copy($ws.Range("H2"))
$ws.Range("B4").Copy($ws.Range("H3"))
$ws.Range("B4").Copy($ws.Range("H4"))
$ws.Range("B4").Copy($ws.Range("H5"))
$ws.Range("B4").Copy($ws.Range("I2"))

$ws.Range("H2").Value = "Built-in Loop-A"
$ws.Range("I2").Value = 2
$ws.Range("H3").Value = "Built-in Loop-B"
$ws.Range("H4").Value = "Built-in Loop-C"
$ws.Range("H5").Value = "Built-in Loop-D"

# --- Update the User Story cell (B4) -----------------------------------
# Clear its formatting (it goes back to the default/no style) and set
# the new text.
$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "NGC-494/T416 OR TC-164"

# --- Update DC Unit Loading Details table data -------------------------
# Row 8 (DDM 800 Loop / Ancillary Conventional)
$ws.Range("D8").Value = 768
$ws.Range("F8").Value = 768

# Row 9 (now FV 411 F / Detectors)
$ws.Range("A9").Value = "FV 411 F"
$ws.Range("B9").Value = "Detectors"
$ws.Range("C9").Value = 117
$ws.Range("D9").Value = 3963
$ws.Range("F9").Value = 3963

# Row 10 (now FV 411 F / Detectors)
$ws.Range("A10").Value = "FV 411 F"
$ws.Range("B10").Value = "Detectors"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 4017
$ws.Range("F10").Value = 4017

# --- Move the selection to B4 ------------------------------------------
$ws.Range("B4").Select() | Out-Null

# --- Size the new "Loop" column to fit its contents ---------------------
$ws.Columns.Item(8).ColumnWidth = 12.5
